$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: update date range, interval, TP/SL %, strategy ---
$ws.Range("D2").Value2 = 44562
$ws.Range("E2").Value2 = 44926
$ws.Range("F2").Value = "3m"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "ScalpEmaRsiAdx"

# --- Row 3: update date range, interval, TP/SL %, strategy ---
$ws.Range("D3").Value2 = 44562
$ws.Range("E3").Value2 = 44926
$ws.Range("F3").Value = "3m"
$ws.Range("G3").Value = 1.2
$ws.Range("H3").Value = 1.2
$ws.Range("I3").Value = "ScalpEmaRsiAdx"

# --- Rows 4-7: revert back to blank template rows (like row 8) ---
$ws.Range("A4:K7").Clear()

$ws.Range("B8").Copy()
$ws.Range("B4:B7").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C4:C7").PasteSpecial(-4122)

$ws.Range("F8").Copy()
$ws.Range("F4:F7").PasteSpecial(-4122)

$ws.Range("I8").Copy()
$ws.Range("I4:I7").PasteSpecial(-4122)

# Rows 4 & 5 used to have an explicit custom height; drop the override so they
# match the other blank template rows.
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

$excel.CutCopyMode = $false

# --- Remove the trailing blank template rows 30-33 ---
$ws.Range("A30:K33").Delete()

# --- Selection / view ---
$ws.Range("E4").Select()
